# Adiciona os novos registros (linhas 9-12) na planilha "Banco de dados"
# conforme atualizacao de 09/11/22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cada linha: nome, fone, sexo(coluna C na planilha mas usada como idade),
# idade(coluna D usada como sexo), endereco -- seguindo exatamente a mesma
# disposicao de colunas ja usada nas linhas existentes da planilha.
$rows = @(
    @("Celso",             "31 985615649",  "40", "Masculino", "Rua : Salvia 162"),
    @("Bruno Marcelino ",  "31 9 85467898", "47", "Femenino",  "Rua :  do pau comeu "),
    @("Breno Josefino ",   "31 9 87658906", "19", "Masculino", "Rua ; Faca na caveira "),
    @("Lindeia",           "234234234",     "12", "Femenino",  "23432423423423")
)

$startRow = 9
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    # Os demais campos da planilha sao sempre texto (mesmo quando parecem
    # numeros, como a idade ou o numero da rua) -- usamos o prefixo de
    # apostrofo para impedir que o Excel os converta em valores numericos.
    $ws.Cells.Item($r, 3).Value = "'" + $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# B12 e E12 tambem sao puramente numericos nesta atualizacao; protege-los
# individualmente para que permanecam como texto, igual ao restante da
# planilha.
$ws.Cells.Item(12, 2).Value = "'234234234"
$ws.Cells.Item(12, 5).Value = "'23432423423423"
